$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.1700000000005
$ws.Range("G2").Value = 0.00012450012692044601
$ws.Range("H2").Value = 0.00069147936187401735
$ws.Range("K2").Value = 5.295015933888039
$ws.Range("L2").Value = "[1.962331476942179, 8.6277003908339]"
$ws.Range("M2").Value = 0.00193199245420649902
$ws.Range("N2").Value = 0.00193199245420649902
$ws.Range("O2").Value = -1.081789662497386
$ws.Range("P2").Value = "[-1.7107371406935412, -0.4528421843012316]"
$ws.Range("Q2").Value = 0.00080121291144563678
$ws.Range("R2").Value = 0.00080121291144563678
$ws.Range("S2").Value = 13.97383856003016
$ws.Range("T2").Value = "[12.25334885463578, 15.69432826542454]"
$ws.Range("W2").Value = 4.333573573573659
$ws.Range("X2").Value = 1.81405405405409
$ws.Range("Y2").Value = 6.853093093093229

# Row 3 updates
$ws.Range("E3").Value = 23.9500000000003
$ws.Range("G3").Value = 0.00004682166564984325
$ws.Range("H3").Value = 0.00069147936187401735
$ws.Range("K3").Value = 4.751125961729254
$ws.Range("L3").Value = "[2.1114245359466786, 7.390827387511829]"
$ws.Range("M3").Value = 0.00044773004165432623
$ws.Range("N3").Value = 0.00089546008330865234
$ws.Range("O3").Value = 2.861711025792505
$ws.Range("P3").Value = "[2.220184598032427, 3.5032374535525825]"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 13.58609779486475
$ws.Range("T3").Value = "[12.123974737949817, 15.048220851779673]"
$ws.Range("W3").Value = 13.04184184184201
$ws.Range("X3").Value = 10.59649649649663
$ws.Range("Y3").Value = 15.48718718718738
